$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 17 (copy down formatting from row 16 first so the
#    inserted row keeps the table's borders/fills/number formats).
$ws.Range("B16:J16").Copy()
$ws.Rows.Item(17).Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Populate the new row 17 with the second period (2509) for the same
#    worker/employer, mirroring row 16's data.
$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "1047393721"
$ws.Range("D17").Value2 = "ERNITH PATERNINA BEDOYA"
$ws.Range("E17").Value2 = "2509"
$ws.Range("F17").Value2 = 56940
$ws.Range("G17").Value2 = 1423500

# 3) Update the summary figures: two overdue periods now, so the total
#    overdue amount doubles and the period counter goes from 1 to 2.
$ws.Range("E11").Value2 = 113880
$ws.Range("F13").Value2 = 2
